$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New lesson-2 sentence pattern rows appended below the existing table.
# Cell values are entered in the same order the author originally typed
# them in (this governs the shared-string table insertion order).
$ws.Range("A10").Value = "これは じしょ です"
$ws.Range("C10").Value = "這是我的辭典"
$ws.Range("A11").Value = "それは わたし の かさ です"
$ws.Range("A12").Value = "この ほん は わたし の です"
$ws.Range("C13").Value = "這是原子筆嗎"
$ws.Range("A9").Value = "*第2課文型"
$ws.Range("C11").Value = "那是我的雨傘"
$ws.Range("B12").Value = "この 本 は 私 の です"
$ws.Range("C12").Value = "這本書是我的"
$ws.Range("B10").Value = "これは 辞書 です"
$ws.Range("B11").Value = "それは 私の傘 です"
$ws.Range("A13").Value = "これは ボルーペンですか"

# Match the style used for the section header (A9) and the last two
# Chinese-translation cells (C11, C12), mirroring the author's formatting.
$ws.Range("A9").Style = $ws.Range("A2").Style
$ws.Range("C11").Style = $ws.Range("A2").Style
$ws.Range("C12").Style = $ws.Range("A2").Style

$ws.Range("B15").Select()
